$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: extend with two new columns (P1=14, Q1=15),
# matching the bold/centered/bordered style used by the rest of row 1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-25: swap the I/K/M/O values (1<->2) and append the new
# P and Q columns (both value 2).
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
